# Adds rows 40-78 (new horimetro entries) to the worksheet and
# extends the used range from A1:L39 to A1:L78.
#
# Every value is written with a leading literal apostrophe so Excel
# stores it as text (matching the source data, where every column,
# including numeric-looking ids/quantities/dates, is t="str").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one new row, columns A..L in order.
$newRows = @(
    ,@('''control system','''17/06/2025, 18:30','''7V3WV','''2025-06-17','''25','''COMPLETO','''100','''5.80','''363','''375','''12.0','''')
    ,@('''control system','''17/06/2025, 18:31','''4YZCP','''2025-06-17','''58','''COMPLETO','''100','''6.54','''250','''265','''15.0','''')
    ,@('''control system','''17/06/2025, 18:31','''Y4FZA','''2025-06-18','''25','''COMPLETO','''100','''5.80','''375','''387','''12.0','''')
    ,@('''control system','''17/06/2025, 18:41','''ZGYT6','''2025-06-17','''25','''COMPLETO','''25','''23.20','''387','''400','''13.0','''')
    ,@('''control system','''17/06/2025, 18:42','''F40RK','''2025-06-10','''48','''COMPLETO','''70','''9.71','''280.2','''303.5','''23.3','''')
    ,@('''control system','''17/06/2025, 18:42','''T6O5H','''2025-06-17','''43','''COMPLETO','''100','''6.79','''0','''15','''15.0','''')
    ,@('''control system','''17/06/2025, 18:45','''B8KS2','''2025-06-17','''25','''COMPLETO','''100','''5.80','''400','''412','''12.0','''')
    ,@('''control system','''17/06/2025, 18:46','''56F5V','''2025-06-18','''25','''COMPLETO','''100','''5.80','''412','''425','''13.0','''')
    ,@('''control system','''17/06/2025, 18:56','''B5Q45','''2025-06-16','''25','''COMPLETO','''100','''5.80','''425','''437','''12.0','''')
    ,@('''control system','''17/06/2025, 19:05','''ESG48','''2025-06-18','''25','''COMPLETO','''100','''5.80','''437','''449','''12.0','''teste de reste inputs')
    ,@('''control system','''17/06/2025, 19:20','''E7NP4','''2025-06-17','''25','''COMPLETO','''100','''5.80','''449','''462','''13.0','''')
    ,@('''control system','''17/06/2025, 19:20','''1IFES','''2025-06-24','''25','''COMPLETO','''100','''5.80','''462','''475','''13.0','''')
    ,@('''control system','''17/06/2025, 19:26','''S23A8','''2025-06-18','''25','''COMPLETO','''100','''5.80','''475','''488','''13.0','''')
    ,@('''control system','''17/06/2025, 19:27','''AHIEI','''2025-06-18','''48','''COMPLETO','''70','''9.71','''303.5','''326','''22.5','''')
    ,@('''control system','''17/06/2025, 19:29','''MCVJL','''2025-06-16','''47','''COMPLETO','''100','''7.42','''0','''15','''15.0','''')
    ,@('''control system','''17/06/2025, 19:29','''P9HE8','''2025-06-24','''70','''COMPLETO','''100','''6.42','''0','''15','''15.0','''')
    ,@('''control system','''17/06/2025, 19:31','''K7E85','''2025-06-18','''25','''COMPLETO','''100','''5.80','''488','''500','''12.0','''')
    ,@('''control system','''18/06/2025, 00:37','''1EHEW','''2025-06-18','''10','''COMPLETO','''100','''6.30','''189','''204','''15.0','''fertirrigação')
    ,@('''control system','''18/06/2025, 01:38','''ZUM2Z','''2025-06-17','''18','''COMPLETO','''100','''6.51','''0','''15','''15.0','''teste de lançamento')
    ,@('''control system','''18/06/2025, 01:40','''46NI7','''2025-06-19','''18','''1/2-PIVO','''70','''9.30','''15','''28','''13.0','''')
    ,@('''control system','''18/06/2025, 01:41','''XCIP7','''2025-06-25','''25','''COMPLETO','''100','''5.80','''500','''515','''15.0','''')
    ,@('''control system','''18/06/2025, 01:42','''HR3UE','''2025-06-19','''25','''COMPLETO','''100','''5.80','''515','''530','''15.0','''')
    ,@('''control system','''18/06/2025, 01:42','''973XR','''2025-06-25','''25','''COMPLETO','''100','''5.80','''530','''545','''15.0','''')
    ,@('''control system','''18/06/2025, 01:43','''E99VT','''2025-06-19','''18','''ABCD','''100','''6.51','''28','''45','''17.0','''')
    ,@('''control system','''18/06/2025, 02:03','''YKN4R','''2025-06-20','''133','''COMPLETO','''50','''10.30','''0','''30','''30.0','''Kálito wilk vieira da silva')
    ,@('''control system','''18/06/2025, 02:05','''VJSAL','''2025-06-25','''133','''COMPLETO','''50','''10.30','''0','''30','''30.0','''segundo teste de lançamento por Kálito Wilk')
    ,@('''control system','''18/06/2025, 02:15','''44LQ9','''2025-06-18','''25','''COMPLETO','''100','''5.80','''545','''560','''15.0','''')
    ,@('''control system','''18/06/2025, 02:15','''6G4FH','''2025-06-17','''25','''COMPLETO','''100','''5.80','''560','''575','''15.0','''')
    ,@('''control system','''18/06/2025, 02:16','''GN18D','''2025-06-25','''25','''COMPLETO','''10','''58.00','''575','''589','''14.0','''')
    ,@('''control system','''18/06/2025, 02:18','''PQ5T5','''2025-06-25','''47','''1/2-PIVO','''100','''7.42','''15','''45','''30.0','''')
    ,@('''control system','''18/06/2025, 02:21','''N2BCM','''2025-06-18','''25','''COMPLETO','''25','''23.20','''589','''600','''11.0','''')
    ,@('''control system','''18/06/2025, 03:01','''NGKBM','''2025-06-18','''25','''COMPLETO','''100','''5.80','''600','''615','''15.0','''')
    ,@('''control system','''18/06/2025, 03:12','''HQTLM','''2025-06-27','''38','''COMPLETO','''100','''8.05','''0','''15','''15.0','''')
    ,@('''control system','''18/06/2025, 03:26','''WW84J','''2025-06-18','''25','''COMPLETO','''100','''5.80','''615','''630','''15.0','''')
    ,@('''control system','''18/06/2025, 03:27','''44SG8','''2025-06-19','''25','''COMPLETO','''100','''5.80','''630','''645','''15.0','''')
    ,@('''control system','''18/06/2025, 03:28','''VC1GH','''2025-06-17','''25','''COMPLETO','''100','''5.80','''645','''660','''15.0','''')
    ,@('''control system','''18/06/2025, 03:29','''1FX6M','''2025-06-11','''15','''COMPLETO','''100','''7.38','''920','''935','''15.0','''')
    ,@('''control system','''18/06/2025, 03:31','''40LD4','''2025-06-28','''15','''COMPLETO','''80','''9.22','''935','''952','''17.0','''')
    ,@('''control system','''18/06/2025, 03:34','''7SDAY','''2025-06-18','''74','''COMPLETO','''100','''6.76','''515','''535','''20.0','''')
)

$startRow = 40
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    $arr = New-Object "object[,]" 1,12
    for ($c = 0; $c -lt 12; $c++) {
        $arr[0, $c] = $rowData[$c]
    }
    $targetRange = $ws.Range($ws.Cells.Item($rowNum, 1), $ws.Cells.Item($rowNum, 12))
    $targetRange.Value = $arr
}
